# repull data, push all data, mean calculation
# Update the dSF column (F) values to reflect newly-repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 6
    4  = 1
    5  = 1
    6  = -4
    7  = -3
    8  = 12
    10 = -1
    12 = -2
    13 = -1
    16 = -1
    17 = 2
    18 = 2
    19 = -1
    20 = -3
    21 = -2
    22 = -1
    23 = -1
    25 = 2
    26 = 1
    27 = 4
    29 = -6
    30 = -4
    31 = 1
    32 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
